# 自动更新Excel文件 - 2025-12-31 23:13:24
# Decrement the "剩余" (remaining days) counter in column E for each row.
# When the remaining-day counter had reached 1 (last day of the cycle), the
# cycle restarts: E resets back to the "总天" (total days) value in column D
# and the start date in column F advances by D days (new cycle begin date).
# All other rows simply get their E value decremented by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose remaining-day count (E) simply decreases by one.
$simpleUpdates = @(
    @{Row=2;  E=11},
    @{Row=3;  E=11},
    @{Row=4;  E=11},
    @{Row=5;  E=3},
    @{Row=6;  E=11},
    @{Row=7;  E=3},
    @{Row=8;  E=11},
    @{Row=9;  E=3},
    @{Row=10; E=4},
    @{Row=11; E=11},
    @{Row=12; E=3},
    @{Row=13; E=11},
    @{Row=14; E=11},
    @{Row=15; E=11},
    @{Row=16; E=7},
    @{Row=17; E=3},
    @{Row=18; E=6},
    @{Row=19; E=6},
    @{Row=20; E=6},
    @{Row=21; E=6},
    @{Row=22; E=3},
    @{Row=23; E=3},
    @{Row=24; E=3},
    @{Row=25; E=3},
    @{Row=26; E=3},
    @{Row=27; E=5},
    @{Row=28; E=6},
    @{Row=29; E=6},
    @{Row=30; E=6},
    @{Row=31; E=6},
    @{Row=32; E=6},
    @{Row=33; E=6},
    @{Row=34; E=6},
    @{Row=35; E=6},
    @{Row=37; E=6},
    @{Row=38; E=6},
    @{Row=39; E=6},
    @{Row=40; E=4},
    @{Row=41; E=4},
    @{Row=42; E=6},
    @{Row=43; E=3},
    @{Row=44; E=4},
    @{Row=45; E=3},
    @{Row=46; E=4},
    @{Row=47; E=6},
    @{Row=48; E=4},
    @{Row=49; E=5},
    @{Row=50; E=1},
    @{Row=51; E=1},
    @{Row=52; E=1},
    @{Row=53; E=1},
    @{Row=54; E=1},
    @{Row=55; E=1},
    @{Row=56; E=1},
    @{Row=57; E=1},
    @{Row=58; E=5},
    @{Row=59; E=5},
    @{Row=60; E=5},
    @{Row=61; E=5},
    @{Row=62; E=5},
    @{Row=63; E=5},
    @{Row=64; E=5},
    @{Row=65; E=6},
    @{Row=66; E=6},
    @{Row=67; E=6},
    @{Row=68; E=6},
    @{Row=69; E=6},
    @{Row=70; E=7},
    @{Row=71; E=7},
    @{Row=72; E=7},
    @{Row=73; E=7},
    @{Row=74; E=7},
    @{Row=75; E=7},
    @{Row=76; E=7},
    @{Row=87; E=4},
    @{Row=88; E=4},
    @{Row=89; E=4},
    @{Row=90; E=4},
    @{Row=91; E=3},
    @{Row=92; E=4},
    @{Row=95; E=9},
    @{Row=96; E=7},
    @{Row=97; E=7},
    @{Row=98; E=7},
    @{Row=99; E=7}
)

# Rows whose cycle rolled over: E resets to the "总天" (D) value and the
# start date (F) jumps forward by D days.
$cycleResets = @(
    @{Row=77; E=10; F=20260101},
    @{Row=78; E=10; F=20260101},
    @{Row=79; E=10; F=20260101},
    @{Row=80; E=10; F=20260101},
    @{Row=81; E=10; F=20260101},
    @{Row=82; E=10; F=20260101},
    @{Row=83; E=10; F=20260101},
    @{Row=84; E=10; F=20260101},
    @{Row=85; E=10; F=20260101},
    @{Row=86; E=10; F=20260101},
    @{Row=93; E=10; F=20260101},
    @{Row=94; E=7;  F=20260101}
)

foreach ($u in $simpleUpdates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

foreach ($u in $cycleResets) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}
